$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date-like text to be stored as a literal
# string (matching the other Date-column cells) instead of being
# auto-parsed into a date serial number by the COM value setter.
$ws.Range("A70").Value = "'2025-10-24"
$ws.Range("A70").ClearFormats()

$ws.Range("B70").Value = 53.81999969482422
$ws.Range("C70").Value = 405.7999877929688
$ws.Range("D70").Value = 326.6000061035156
